# Commit: "update data settings to activate ptdf calculation"
#
# This updates the "rel_node__commodity" sheet: most node rows get their
# commodity changed from "electricity" to a new "energy_carrier" value
# (rows 2 and 6-86; rows 3-5, the region nodes, are left as "electricity").
# It also switches the active/selected worksheet from "obj_connection" to
# "rel_node__commodity", with a new selection on that sheet.

$wb = $excel.ActiveWorkbook

$wsCommodity = $wb.Worksheets.Item("rel_node__commodity")

# Row 2 plus rows 6 through 86 (column C = commodity) move from
# "electricity" (shared string 241) to the newly introduced
# "energy_carrier" shared string. Rows 3-5 (regions) stay untouched.
$wsCommodity.Cells.Item(2, 3).Value = "energy_carrier"
for ($r = 6; $r -le 86; $r++) {
    $wsCommodity.Cells.Item($r, 3).Value = "energy_carrier"
}

# Activate the sheet and move the selection, mirroring the saved view
# state in the workbook (this also flips tabSelected from
# "obj_connection" onto "rel_node__commodity" and updates the workbook's
# activeTab).
$wsCommodity.Activate()
$wsCommodity.Range("E6").Select()
